$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the is_locked_lbl / is_enabled_lbl template cells (D1:E1),
# shifting the remaining cells (order_by, rem) left.
$ws.Range("D1:E1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
